# Scheduled runner update: refresh market-board price/profit figures
# (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ]) on
# the Leve-profit sheets, per latest Universalis data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 3945.2
$ws.Range("I9").Value = 4638.3477
$ws.Range("J9").Value = 2616.6667
$ws.Range("K9").Value = 4638.3477
$ws.Range("L9").Value = 2616.6667
$ws.Range("M9").Value = -4469.3477
$ws.Range("N9").Value = -2954.6667
$ws.Range("H100").Value = 2257.111
$ws.Range("I100").Value = 2120.8333
$ws.Range("J100").Value = 2529.6667
$ws.Range("K100").Value = 2120.8333
$ws.Range("L100").Value = 2529.6667
$ws.Range("M100").Value = -1579.8333
$ws.Range("N100").Value = -3611.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2531.889
$ws.Range("I2").Value = 1993.75
$ws.Range("K2").Value = 1993.75
$ws.Range("M2").Value = -1880.75
$ws.Range("H61").Value = 76924640
$ws.Range("I61").Value = 90910470
$ws.Range("J61").Value = 2547.5
$ws.Range("K61").Value = 90910470
$ws.Range("L61").Value = 2547.5
$ws.Range("M61").Value = -90910258
$ws.Range("N61").Value = -2971.5
$ws.Range("H74").Value = 17545824
$ws.Range("I74").Value = 40001164
$ws.Range("K74").Value = 40001164
$ws.Range("M74").Value = -40000290
$ws.Range("H77").Value = 17545824
$ws.Range("I77").Value = 40001164
$ws.Range("K77").Value = 200005820
$ws.Range("M77").Value = -200001452
$ws.Range("H80").Value = 59945
$ws.Range("J80").Value = 59945
$ws.Range("L80").Value = 59945
$ws.Range("N80").Value = -61941
$ws.Range("H83").Value = 59945
$ws.Range("J83").Value = 59945
$ws.Range("L83").Value = 179835
$ws.Range("N83").Value = -189819
$ws.Range("H116").Value = 2531.889
$ws.Range("I116").Value = 1993.75
$ws.Range("K116").Value = 1993.75
$ws.Range("M116").Value = 300.25
$ws.Range("H136").Value = 76924640
$ws.Range("I136").Value = 90910470
$ws.Range("J136").Value = 2547.5
$ws.Range("K136").Value = 272731410
$ws.Range("L136").Value = 7642.5
$ws.Range("M136").Value = -272728860
$ws.Range("N136").Value = -12742.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2531.889
$ws.Range("I3").Value = 1993.75
$ws.Range("K3").Value = 1993.75
$ws.Range("M3").Value = -1879.75
$ws.Range("H82").Value = 41317
$ws.Range("J82").Value = 49189.332
$ws.Range("L82").Value = 49189.332
$ws.Range("N82").Value = -49955.332
$ws.Range("H85").Value = 41317
$ws.Range("J85").Value = 49189.332
$ws.Range("L85").Value = 49189.332
$ws.Range("N85").Value = -51841.332
$ws.Range("H86").Value = 7562.2856
$ws.Range("I86").Value = 8787.9375
$ws.Range("K86").Value = 8787.9375
$ws.Range("M86").Value = -7664.9375
$ws.Range("H89").Value = 7562.2856
$ws.Range("I89").Value = 8787.9375
$ws.Range("K89").Value = 43939.6875
$ws.Range("M89").Value = -38323.6875
$ws.Range("H105").Value = 10213.105
$ws.Range("I105").Value = 22751.5
$ws.Range("K105").Value = 22751.5
$ws.Range("M105").Value = -21004.5
$ws.Range("H107").Value = 3304
$ws.Range("I107").Value = 2099.077
$ws.Range("K107").Value = 2099.077
$ws.Range("M107").Value = -179.0770000000002
$ws.Range("H134").Value = 1931.4222
$ws.Range("I134").Value = 2204.8235
$ws.Range("K134").Value = 6614.470499999999
$ws.Range("M134").Value = -4079.470499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 41671304
$ws.Range("I31").Value = 2676.1428
$ws.Range("J31").Value = 138898110
$ws.Range("K31").Value = 2676.1428
$ws.Range("L31").Value = 138898110
$ws.Range("M31").Value = -2381.1428
$ws.Range("N31").Value = -138898700
$ws.Range("H34").Value = 41671304
$ws.Range("I34").Value = 2676.1428
$ws.Range("J34").Value = 138898110
$ws.Range("K34").Value = 2676.1428
$ws.Range("L34").Value = 138898110
$ws.Range("M34").Value = -2474.1428
$ws.Range("N34").Value = -138898514
$ws.Range("H140").Value = 63244
$ws.Range("J140").Value = 95779
$ws.Range("L140").Value = 95779
$ws.Range("N140").Value = -106139
$ws.Range("H141").Value = 85049.94500000001
$ws.Range("J141").Value = 85049.94500000001
$ws.Range("L141").Value = 85049.94500000001
$ws.Range("N141").Value = -95409.94500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 14846.137
$ws.Range("I56").Value = 14846.137
$ws.Range("K56").Value = 14846.137
$ws.Range("M56").Value = -14316.137
$ws.Range("H63").Value = 1000
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 1000
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H128").Value = 188659.33
$ws.Range("I128").Value = 188659.33
$ws.Range("K128").Value = 565977.99
$ws.Range("M128").Value = -560997.99

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 15504.286
$ws.Range("I29").Value = 12000
$ws.Range("J29").Value = 16088.333
$ws.Range("K29").Value = 12000
$ws.Range("L29").Value = 16088.333
$ws.Range("M29").Value = -11710
$ws.Range("N29").Value = -16668.333
$ws.Range("H58").Value = 46000.5
$ws.Range("J58").Value = 48667.332
$ws.Range("L58").Value = 48667.332
$ws.Range("N58").Value = -49221.332
$ws.Range("H107").Value = 549.8077
$ws.Range("J107").Value = 404.5
$ws.Range("L107").Value = 404.5
$ws.Range("N107").Value = -4244.5
$ws.Range("H122").Value = 38463984
$ws.Range("I122").Value = 2073.0908
$ws.Range("K122").Value = 6219.2724
$ws.Range("M122").Value = -3769.2724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 628.25
$ws.Range("I55").Value = 482.6
$ws.Range("K55").Value = 482.6
$ws.Range("M55").Value = -309.6
$ws.Range("H68").Value = 2685.4546
$ws.Range("I68").Value = 1985.625
$ws.Range("J68").Value = 4551.6665
$ws.Range("K68").Value = 1985.625
$ws.Range("L68").Value = 4551.6665
$ws.Range("M68").Value = -1236.625
$ws.Range("N68").Value = -6049.6665
$ws.Range("H71").Value = 2685.4546
$ws.Range("I71").Value = 1985.625
$ws.Range("J71").Value = 4551.6665
$ws.Range("K71").Value = 9928.125
$ws.Range("L71").Value = 22758.3325
$ws.Range("M71").Value = -6184.125
$ws.Range("N71").Value = -30246.3325
$ws.Range("H82").Value = 4250.75
$ws.Range("J82").Value = 3667.6667
$ws.Range("L82").Value = 3667.6667
$ws.Range("N82").Value = -4389.6667
$ws.Range("H85").Value = 4250.75
$ws.Range("J85").Value = 3667.6667
$ws.Range("L85").Value = 3667.6667
$ws.Range("N85").Value = -6163.6667
$ws.Range("H131").Value = 86174.336
$ws.Range("J131").Value = 102937.5
$ws.Range("L131").Value = 102937.5
$ws.Range("N131").Value = -113017.5
$ws.Range("H139").Value = 64872
$ws.Range("J139").Value = 69419.5
$ws.Range("L139").Value = 69419.5
$ws.Range("N139").Value = -79699.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1825000
$ws.Range("J2").Value = 50000
$ws.Range("L2").Value = 50000
$ws.Range("N2").Value = -50224
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H126").Value = 8692.666999999999
$ws.Range("I126").Value = 8527.929
$ws.Range("K126").Value = 25583.787
$ws.Range("M126").Value = -23113.787
$ws.Range("H132").Value = 1568.6428
$ws.Range("I132").Value = 1595.9
$ws.Range("K132").Value = 4787.700000000001
$ws.Range("M132").Value = -2257.700000000001
